$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.944
$ws.Range("A21").Value = -19.884
$ws.Range("A23").Value = -20.317
$ws.Range("A25").Value = -21.734
$ws.Range("E27").Value = 16.391
$ws.Range("E31").Value = 16.62
$ws.Range("E39").Value = 16.401
$ws.Range("E48").Value = 17.062
$ws.Range("E51").Value = 16.617
$ws.Range("E52").Value = 16.543
$ws.Range("A53").Value = -22.013
$ws.Range("E55").Value = 16.416
$ws.Range("E56").Value = 16.276
$ws.Range("A57").Value = -22.17
$ws.Range("E57").Value = 16.453
$ws.Range("A59").Value = -22.5
$ws.Range("A69").Value = -21.648
$ws.Range("E73").Value = 16.572
$ws.Range("A79").Value = -21.23
$ws.Range("A83").Value = -21.997
$ws.Range("E89").Value = 17.362
$ws.Range("E90").Value = 16.537
$ws.Range("A93").Value = -21.524
